$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting MSSV/Lời nhắn to the right
$ws.Range("B1").EntireColumn.Insert()

# Set the header for the newly inserted column
$ws.Range("B1").Value = "Ngày"

# Update the active selection to B1 (single cell) as in the edited workbook
$ws.Range("B1").Select()
